$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.695.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.200.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '257.24'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '84.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +13.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.615'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.595'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0920'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.98%  '

$ws.Range("E13").Value = '  +2.20%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.530.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("E15").Value = '  +1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.242.18'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.47%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.784'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.627.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000103'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.44%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.86%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.38'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0862'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.31'
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = '  +1.82%  '

$ws.Range("E36").Value = '  +2.18%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0359'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.58%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.48'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.66%  '

$ws.Range("E41").Value = '  +0.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '63.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.50%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.45%  '

$ws.Range("E44").Value = '  +1.42%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.99%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0981'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.58%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("E48").Value = '  +4.81%  '

$ws.Range("E49").Value = '  +0.99%  '

$ws.Range("E50").Value = '  -5.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.77%  '
